$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 5032.778
$ws.Range("J58").Value = 6264.048
$ws.Range("L58").Value = 18792.144
$ws.Range("N58").Value = -19092.144

$ws.Range("H64").Value = 4722.0557
$ws.Range("I64").Value = 3499
$ws.Range("J64").Value = 4966.6665
$ws.Range("K64").Value = 3499
$ws.Range("L64").Value = 4966.6665
$ws.Range("M64").Value = -3251
$ws.Range("N64").Value = -5462.6665

$ws.Range("H67").Value = 4722.0557
$ws.Range("I67").Value = 3499
$ws.Range("J67").Value = 4966.6665
$ws.Range("K67").Value = 3499
$ws.Range("L67").Value = 4966.6665
$ws.Range("M67").Value = -2641
$ws.Range("N67").Value = -6682.6665

$ws.Range("H88").Value = 436980.62
$ws.Range("I88").Value = 40299.6
$ws.Range("J88").Value = 578652.4399999999
$ws.Range("K88").Value = 40299.6
$ws.Range("L88").Value = 578652.4399999999
$ws.Range("M88").Value = -39893.6
$ws.Range("N88").Value = -579464.4399999999

$ws.Range("H91").Value = 436980.62
$ws.Range("I91").Value = 40299.6
$ws.Range("J91").Value = 578652.4399999999
$ws.Range("K91").Value = 40299.6
$ws.Range("L91").Value = 578652.4399999999
$ws.Range("M91").Value = -38895.6
$ws.Range("N91").Value = -581460.4399999999

$ws.Range("H96").Value = 1926.0769
$ws.Range("I96").Value = 1060.2222
$ws.Range("J96").Value = 3874.25
$ws.Range("K96").Value = 3180.6666
$ws.Range("L96").Value = 11622.75
$ws.Range("M96").Value = -1807.6666
$ws.Range("N96").Value = -14368.75

$ws.Range("H137").Value = 299238.1
$ws.Range("I137").Value = 669457.1
$ws.Range("J137").Value = 6959.8945
$ws.Range("K137").Value = 2008371.3
$ws.Range("L137").Value = 20879.6835
$ws.Range("M137").Value = -2005821.3
$ws.Range("N137").Value = -25979.6835

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 23812002
$ws.Range("I45").Value = 26317738
$ws.Range("J45").Value = 7507
$ws.Range("K45").Value = 26317738
$ws.Range("L45").Value = 7507
$ws.Range("M45").Value = -26317361
$ws.Range("N45").Value = -8261

$ws.Range("H61").Value = 7815674.5
$ws.Range("I61").Value = 3289.4333
$ws.Range("K61").Value = 3289.4333
$ws.Range("M61").Value = -3077.4333

$ws.Range("H74").Value = 7411.4736
$ws.Range("I74").Value = 2920.7036
$ws.Range("K74").Value = 2920.7036
$ws.Range("M74").Value = -2046.7036

$ws.Range("H77").Value = 7411.4736
$ws.Range("I77").Value = 2920.7036
$ws.Range("K77").Value = 14603.518
$ws.Range("M77").Value = -10235.518

$ws.Range("H110").Value = 1805.6666
$ws.Range("I110").Value = 1787.625
$ws.Range("J110").Value = 1950
$ws.Range("K110").Value = 1787.625
$ws.Range("L110").Value = 1950
$ws.Range("M110").Value = 257.375
$ws.Range("N110").Value = -6040

$ws.Range("H132").Value = 5716.1377
$ws.Range("I132").Value = 3447.8635
$ws.Range("J132").Value = 12845
$ws.Range("K132").Value = 10343.5905
$ws.Range("L132").Value = 38535
$ws.Range("M132").Value = -7813.5905
$ws.Range("N132").Value = -43595

$ws.Range("H136").Value = 7815674.5
$ws.Range("I136").Value = 3289.4333
$ws.Range("K136").Value = 9868.2999
$ws.Range("M136").Value = -7318.2999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 3264.4546
$ws.Range("I22").Value = 2190.6667
$ws.Range("J22").Value = 4553
$ws.Range("K22").Value = 2190.6667
$ws.Range("L22").Value = 4553
$ws.Range("M22").Value = -2017.6667
$ws.Range("N22").Value = -4899

$ws.Range("H35").Value = 57500
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

$ws.Range("H105").Value = 2314.8333
$ws.Range("J105").Value = 3000
$ws.Range("L105").Value = 3000
$ws.Range("N105").Value = -6494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 956.6667
$ws.Range("I16").Value = 968.4
$ws.Range("J16").Value = 898
$ws.Range("K16").Value = 968.4
$ws.Range("L16").Value = 898
$ws.Range("M16").Value = -681.4
$ws.Range("N16").Value = -1472

$ws.Range("H31").Value = 870653.7
$ws.Range("I31").Value = 10941.444
$ws.Range("J31").Value = 2590078
$ws.Range("K31").Value = 10941.444
$ws.Range("L31").Value = 2590078
$ws.Range("M31").Value = -10646.444
$ws.Range("N31").Value = -2590668

$ws.Range("H34").Value = 870653.7
$ws.Range("I34").Value = 10941.444
$ws.Range("J34").Value = 2590078
$ws.Range("K34").Value = 10941.444
$ws.Range("L34").Value = 2590078
$ws.Range("M34").Value = -10739.444
$ws.Range("N34").Value = -2590482

$ws.Range("H113").Value = 956.6667
$ws.Range("I113").Value = 968.4
$ws.Range("J113").Value = 898
$ws.Range("K113").Value = 968.4
$ws.Range("L113").Value = 898
$ws.Range("M113").Value = 1201.6
$ws.Range("N113").Value = -5238

$ws.Range("H124").Value = 67502
$ws.Range("J124").Value = 67502
$ws.Range("L124").Value = 67502
$ws.Range("N124").Value = -72412

$ws.Range("H134").Value = 2253.5386
$ws.Range("I134").Value = 1135.5938
$ws.Range("K134").Value = 3406.7814
$ws.Range("M134").Value = -871.7814000000003

$ws.Range("H141").Value = 264199.8
$ws.Range("J141").Value = 285406.72
$ws.Range("L141").Value = 285406.72
$ws.Range("N141").Value = -295766.72

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 290
$ws.Range("I49").Value = 290
$ws.Range("K49").Value = 870
$ws.Range("M49").Value = -714

$ws.Range("H86").Value = 551.8570999999999
$ws.Range("I86").Value = 510.75
$ws.Range("J86").Value = 606.6667
$ws.Range("K86").Value = 1532.25
$ws.Range("L86").Value = 1820.0001
$ws.Range("M86").Value = -346.25
$ws.Range("N86").Value = -4192.0001

$ws.Range("H89").Value = 551.8570999999999
$ws.Range("I89").Value = 510.75
$ws.Range("J89").Value = 606.6667
$ws.Range("K89").Value = 4596.75
$ws.Range("L89").Value = 5460.0003
$ws.Range("M89").Value = 1331.25
$ws.Range("N89").Value = -17316.0003

$ws.Range("H131").Value = 16129.167
$ws.Range("I131").Value = 10510
$ws.Range("J131").Value = 20142.857
$ws.Range("K131").Value = 31530
$ws.Range("L131").Value = 60428.571
$ws.Range("M131").Value = -26490
$ws.Range("N131").Value = -70508.571

$ws.Range("H140").Value = 204521.6
$ws.Range("J140").Value = 12000
$ws.Range("L140").Value = 36000
$ws.Range("N140").Value = -46360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4803.9375
$ws.Range("I80").Value = 3433.7144
$ws.Range("J80").Value = 5869.6665
$ws.Range("K80").Value = 3433.7144
$ws.Range("L80").Value = 5869.6665
$ws.Range("M80").Value = -2435.7144
$ws.Range("N80").Value = -7865.6665

$ws.Range("H83").Value = 4803.9375
$ws.Range("I83").Value = 3433.7144
$ws.Range("J83").Value = 5869.6665
$ws.Range("K83").Value = 17168.572
$ws.Range("L83").Value = 29348.3325
$ws.Range("M83").Value = -12176.572
$ws.Range("N83").Value = -39332.3325

$ws.Range("H113").Value = 4380
$ws.Range("I113").Value = 3722.3333
$ws.Range("J113").Value = 5037.6665
$ws.Range("K113").Value = 3722.3333
$ws.Range("L113").Value = 5037.6665
$ws.Range("M113").Value = -1552.3333
$ws.Range("N113").Value = -9377.666499999999

$ws.Range("H132").Value = 3406.8975
$ws.Range("I132").Value = 3293.0625
$ws.Range("K132").Value = 9879.1875
$ws.Range("M132").Value = -7349.1875

$ws.Range("H136").Value = 30000
$ws.Range("J136").Value = 30000
$ws.Range("L136").Value = 90000
$ws.Range("N136").Value = -95100

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1048.2
$ws.Range("I22").Value = 1122
$ws.Range("K22").Value = 1122
$ws.Range("M22").Value = -827

$ws.Range("H27").Value = 1048.2
$ws.Range("I27").Value = 1122
$ws.Range("K27").Value = 1122
$ws.Range("M27").Value = -1015

$ws.Range("H46").Value = 4123.769
$ws.Range("I46").Value = 2517.5
$ws.Range("K46").Value = 2517.5
$ws.Range("M46").Value = -2329.5

$ws.Range("H93").Value = 1986.1154
$ws.Range("J93").Value = 2480.3333
$ws.Range("L93").Value = 2480.3333
$ws.Range("N93").Value = -4976.3333

$ws.Range("H132").Value = 6473.4287
$ws.Range("I132").Value = 9284.5
$ws.Range("J132").Value = 2725.3333
$ws.Range("K132").Value = 27853.5
$ws.Range("L132").Value = 8175.999899999999
$ws.Range("M132").Value = -25323.5
$ws.Range("N132").Value = -13235.9999

$ws.Range("H136").Value = 50796.348
$ws.Range("I136").Value = 11564.25
$ws.Range("J136").Value = 84423.86
$ws.Range("K136").Value = 34692.75
$ws.Range("L136").Value = 253271.58
$ws.Range("M136").Value = -32142.75
$ws.Range("N136").Value = -258371.58

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 39116.45
$ws.Range("J54").Value = 34247.5
$ws.Range("L54").Value = 34247.5
$ws.Range("N54").Value = -35287.5

$ws.Range("H100").Value = 795.3461
$ws.Range("I100").Value = 773.087
$ws.Range("K100").Value = 1546.174
$ws.Range("M100").Value = -1005.174

$ws.Range("H107").Value = 1099.258
$ws.Range("I107").Value = 1301.8572
$ws.Range("J107").Value = 673.8
$ws.Range("K107").Value = 3905.5716
$ws.Range("L107").Value = 2021.4
$ws.Range("M107").Value = -1985.5716
$ws.Range("N107").Value = -5861.4

$ws.Range("H113").Value = 598.3333
$ws.Range("J113").Value = 598.3333
$ws.Range("L113").Value = 1794.9999
$ws.Range("N113").Value = -6134.9999

$ws.Range("H114").Value = 112000
$ws.Range("J114").Value = 112000
$ws.Range("L114").Value = 112000
$ws.Range("N114").Value = -120678

$ws.Range("H132").Value = 2455312.2
$ws.Range("I132").Value = 4036.5518
$ws.Range("K132").Value = 12109.6554
$ws.Range("M132").Value = -9579.6554
